$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column P (rows 4-14) into column Q so the new
# "2020" column inherits the same per-row styles (header style, normal
# data style, bold/bottom-border total-row style, etc.) that the
# workbook author already uses for the other year columns.
$ws.Range("P4:P14").Copy() | Out-Null
$ws.Range("Q4:Q14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New "2020" year header.
$ws.Range("Q4").Value = 2020

# New data values for 2020, one per region/row.
$ws.Range("Q5").Value = 109.7221295941265
$ws.Range("Q6").Value = 108.44905375816947
$ws.Range("Q7").Value = 109.90982951756889
$ws.Range("Q8").Value = 108.40606487500015
$ws.Range("Q9").Value = 109.40161876466024
$ws.Range("Q10").Value = 107.71155656686271
$ws.Range("Q11").Value = 111.78921596090774
$ws.Range("Q12").Value = 111.39254046803097
$ws.Range("Q13").Value = 110.44919152842827
$ws.Range("Q14").Value = 106.89826464456031

# Matches the saved cursor position recorded in the diff.
$ws.Range("N14").Select()
